{"js": "// Edit 1: remove \"Monte Carlo Sampling, \" from the research-interest sentence\n// (\"...Bayesian Learning, Monte Carlo Sampling, Tracking, ...\" ->\n//  \"...Bayesian Learning, Tracking, ...\")\nconst body = context.document.body;\n\nconst hits1 = body.search(\"Bayesian Learning, Monte Carlo Sampling, Tracking,\", { matchCase: true });\nhits1.load(\"text\");\nawait context.sync();\n\nif (hits1.items.length > 0) {\n  hits1.items[0].insertText(\"Bayesian Learning, Tracking,\", \"Replace\");\n  await context.sync();\n}\n\n// The document's \"_GoBack\" bookmark (Word's \"last edit location\" marker) sat\n// right before \"Reinforcement Learning\"; since the edit above now happens\n// earlier in the sentence, Word re-stamps \"_GoBack\" at that newer edit spot \u2014\n// right after \"Bayesian Learning, \" (immediately before \"Tracking,\").\nconst anchorHits = body.search(\"Bayesian Learning, \", { matchCase: true });\nanchorHits.load(\"text\");\nawait context.sync();\n\nif (anchorHits.items.length > 0) {\n  const caret = anchorHits.items[0].getRange(\"End\");\n  context.document.deleteBookmark(\"_GoBack\");\n  caret.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// Edit 2: expand \"Deep Reinforcement Learning, \" (the red \"looking for a job\n// related to ...\" sentence) into \"Deep Reinforcement Learning (Robot Motion\n// Planning), \"\nconst hits2 = body.search(\"Deep Reinforcement Learning, \", { matchCase: true });\nhits2.load(\"text\");\nawait context.sync();\n\nif (hits2.items.length > 0) {\n  hits2.items[0].insertText(\"Deep Reinforcement Learning (Robot Motion Planning), \", \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Edit 1: remove \"Monte Carlo Sampling, \" from the research-interest sentence\n# (\"...Bayesian Learning, Monte Carlo Sampling, Tracking, ...\" ->\n#  \"...Bayesian Learning, Tracking, ...\")\n$find1 = $d.Content\n$find1.Find.Execute(\"Bayesian Learning, Monte Carlo Sampling, Tracking,\", $false, $false, $false, $false, $false, $true, 1, $false, \"Bayesian Learning, Tracking,\", 2)\n\n# The document's \"_GoBack\" bookmark (Word's \"last edit location\" marker) sat\n# right before \"Reinforcement Learning\"; since the edit above now happens\n# earlier in the sentence, Word re-stamps \"_GoBack\" at that newer edit spot -\n# right after \"Bayesian Learning, \" (immediately before \"Tracking,\").\n$anchor = $d.Content\n$anchor.Find.Execute(\"Bayesian Learning, \")\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n$caret = $d.Range($anchor.End, $anchor.End)\n$d.Bookmarks.Add(\"_GoBack\", $caret)\n\n# Edit 2: expand \"Deep Reinforcement Learning, \" (the red \"looking for a job\n# related to ...\" sentence) into \"Deep Reinforcement Learning (Robot Motion\n# Planning), \"\n$find2 = $d.Content\n$find2.Find.Execute(\"Deep Reinforcement Learning, \", $false, $false, $false, $false, $false, $true, 1, $false, \"Deep Reinforcement Learning (Robot Motion Planning), \", 2)\n\nWrite-Output \"done\"\n"}
